$d = $word.ActiveDocument

# --- Edit 1: merge the two runs "et grâce à la propriété sin" + "(arctan) :"
# into a single run, removing the _GoBack bookmark that sat between them.
# A literal Find/Replace across the run boundary naturally merges the runs
# and drops the (now interior) bookmark.
$d.Content.Find.Execute(
    "et grâce à la propriété sin(arctan) :",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "et grâce à la propriété sin(arctan) :",
    2) | Out-Null

# --- Edit 2: append a new italic, 9pt paragraph citing the algorithm source,
# right after the final paragraph ("... à partir des données de départ."),
# and re-insert the _GoBack bookmark at the end of that (now prior) paragraph.
$last = $d.Paragraphs.Last
$lastRange = $last.Range
$lastRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $lastRange) | Out-Null

$newPara = $lastRange.InsertParagraphAfter()
$lastRange.Collapse(0)

$srcRange = $d.Paragraphs.Last.Range
$srcRange.Text = "Source de l’algorithme : https://github.com/cboulay/PSMoveService/wiki/Optical-Tracker-Algorithms"
$srcRange.Font.Italic = $true
$srcRange.Font.Size = 9
